$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(14495006, "2025-09-01", "Naomi Osaka", "Coco Gauff", "Gana Coco Gauff", 1.62),
    @(14581495, "2025-09-01", "Jelle Sels", "Diego Dedura-Palomero", "Gana Jelle Sels", 2.75),
    @(14581491, "2025-09-02", "Andrew Paulson", "Jacopo Berrettini", "Gana Andrew Paulson", 3),
    @(14581504, "2025-09-02", "Joel Schwaerzler", "Tiago Pereira", "Gana Tiago Pereira", 2.63),
    @(14579776, "2025-09-02", "Ilya Ivashka", "Nicolas Mejia", "Gana Ilya Ivashka", 1.83),
    @(14579770, "2025-09-02", "Marc-Andrea Huesler", "Alastair Gray", "Gana Alastair Gray", 3.5),
    @(14579773, "2025-09-02", "Marek Gengel", "Benjamin Hassan", "Gana Marek Gengel", 3.4),
    @(14583803, "2025-09-02", "Rodrigo Pacheco Mendez", "Norbert Gombos", "Gana Norbert Gombos", 3)
)

$startRow = 55
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    # keep 'fecha' as literal text (YYYY-MM-DD), matching the rest of the sheet
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]

    # 'resultado' / 'profit' stay blank (pending result), same as other
    # not-yet-settled rows. Touch/reset the border so the cells are
    # materialised (present, default style) instead of being omitted.
    $ws.Cells.Item($row, 7).Borders.LineStyle = 1
    $ws.Cells.Item($row, 7).Borders.LineStyle = -4142
    $ws.Cells.Item($row, 8).Borders.LineStyle = 1
    $ws.Cells.Item($row, 8).Borders.LineStyle = -4142
}
